# Applies:
#   1) Removal of the "Heading2" paragraph style from the section heading
#      paragraphs (Introduction, Identifying Inequalities..., Policies
#      Addressing Educational Inequalities, Sociological Theories
#      Application, Conclusion) so they fall back to the default
#      (un-styled) paragraph, matching the target OOXML which has no
#      <w:pPr> on these paragraphs at all.
#   2) Citation-text swaps inside specific body paragraphs (the same
#      author-name citation string maps to different replacement
#      citations depending on which paragraph it appears in, so each
#      replacement is scoped to that paragraph's Range rather than done
#      as a blanket document-wide Find/Replace).

$d = $word.ActiveDocument

function Remove-HeadingStyle($paraIndex) {
    $p = $d.Paragraphs.Item($paraIndex)
    # Range.Text includes the trailing paragraph-mark (CR, char 13) -
    # strip it so it doesn't leak into the <w:t> content.
    $text = $p.Range.Text.TrimEnd([char]13)
    # Re-insert the paragraph as plain WordOpenXML with no <w:pPr>, which
    # drops the Heading2 pStyle while preserving the run text exactly
    # (only flag xml:space="preserve" when the text actually has leading
    # or trailing whitespace that needs protecting, matching how Word
    # itself emits these runs).
    $escaped = $text.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    if ($text -ne $text.Trim()) {
        $tOpen = '<w:t xml:space="preserve">'
    } else {
        $tOpen = '<w:t>'
    }
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r>' + $tOpen + $escaped + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $p.Range.InsertXML($xml)
}

function Replace-InParagraph($paraIndex, $find, $replace) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $r.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

# --- 1) Strip Heading2 styling from the five section headings ---
Remove-HeadingStyle 2   # "Introduction"
Remove-HeadingStyle 6   # "Identifying Inequalities in the U.S. Education System"
Remove-HeadingStyle 12  # "Policies Addressing Educational Inequalities  "
Remove-HeadingStyle 16  # "Sociological Theories Application  "
Remove-HeadingStyle 25  # "Conclusion  "

# --- 2) Citation updates, scoped per paragraph ---

# Paragraph 8: "The U.S. education system is riddled with..."
Replace-InParagraph 8 "(Lewis and Diamond)" "(Ref-s309149)"
Replace-InParagraph 8 "(Dupree and Boykin)" "(Ref-s309149)"

# Paragraph 10: "Furthermore, economic factors play a crucial role..."
Replace-InParagraph 10 "(Hardy and Logan)" "(Nguyen, 2015)"

# Paragraph 14: "In response to entrenched educational disparities..."
Replace-InParagraph 14 "(Diem and Welton)" "(Brown and Garcia)"
Replace-InParagraph 14 "(White et al.)" "(Brown and Garcia)"

# Paragraph 19: "Applying Functionalism to educational inequalities..."
Replace-InParagraph 19 "(Lewis and Diamond)" "(Brown and Garcia)"

# Paragraph 23: "Feminist Theory offers a critical lens..."
Replace-InParagraph 23 "(Diem and Welton)" "(Ref-s395099)"
Replace-InParagraph 23 "(Dupree and Boykin)" "(Ref-s395099)"

Write-Host "Edit complete."
